$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.274.96'
$ws.Range("E2").Value = '  +5.23%  '
$ws.Range("D3").Value = '2.475.13'
$ws.Range("E3").Value = '  +3.33%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = "'486.23"
$ws.Range("E5").Value = '  +5.68%  '
$ws.Range("D6").Value = "'146.36"
$ws.Range("E6").Value = '  +12.39%  '
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = "'0.509"
$ws.Range("E8").Value = '  +5.70%  '
$ws.Range("D9").Value = '2.479.55'
$ws.Range("E9").Value = '  +2.78%  '
$ws.Range("E10").Value = '  +8.72%  '
$ws.Range("D11").Value = "'0.0964"
$ws.Range("E11").Value = '  +2.12%  '
$ws.Range("E12").Value = '  +5.82%  '
$ws.Range("E13").Value = '  +1.50%  '
$ws.Range("D14").Value = '2.902.71'
$ws.Range("E14").Value = '  +2.82%  '
$ws.Range("D15").Value = '56.287.13'
$ws.Range("E15").Value = '  +5.16%  '
$ws.Range("D16").Value = "'21.06"
$ws.Range("E16").Value = '  +7.74%  '
$ws.Range("E17").Value = '  +3.16%  '
$ws.Range("D18").Value = '2.485.99'
$ws.Range("E18").Value = '  +2.53%  '
$ws.Range("D19").Value = "'4.51"
$ws.Range("E19").Value = '  +9.01%  '
$ws.Range("D20").Value = "'10.00"
$ws.Range("E20").Value = '  +7.22%  '
$ws.Range("D21").Value = "'317.15"
$ws.Range("E21").Value = '  +3.86%  '
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").Value = "'5.78"
$ws.Range("E23").Value = '  +8.74%  '
$ws.Range("D24").Value = "'58.27"
$ws.Range("E24").Value = '  +4.81%  '
$ws.Range("E25").Value = '  +6.97%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("D27").Value = "'0.162"
$ws.Range("E27").Value = '  +5.93%  '
$ws.Range("D28").Value = '2.585.74'
$ws.Range("E28").Value = '  +3.76%  '
$ws.Range("E29").Value = '  +7.50%  '
$ws.Range("D30").Value = '0.0₃0787'
$ws.Range("E30").Value = '  +10.63%  '
$ws.Range("E31").Value = '  +0.18%  '
$ws.Range("D32").Value = "'148.96"
$ws.Range("E32").Value = '  +2.54%  '
$ws.Range("D33").Value = "'18.09"
$ws.Range("E33").Value = '  +2.95%  '
$ws.Range("E34").Value = '  +6.03%  '
$ws.Range("D35").Value = "'5.18"
$ws.Range("E35").Value = '  +4.51%  '
$ws.Range("E36").Value = '  +8.61%  '
$ws.Range("D37").Value = "'3.71"
$ws.Range("E37").Value = '  +6.07%  '
$ws.Range("D38").Value = "'0.859"
$ws.Range("E38").Value = '  +8.09%  '
$ws.Range("D39").Value = "'34.14"
$ws.Range("E39").Value = '  +4.59%  '
$ws.Range("D40").Value = "'3.50"
$ws.Range("E40").Value = '  +8.06%  '
$ws.Range("D41").Value = "'0.0556"
$ws.Range("E41").Value = '  +6.73%  '
$ws.Range("D42").Value = "'0.994"
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("D43").Value = "'0.607"
$ws.Range("E43").Value = '  +2.81%  '
$ws.Range("E44").Value = '  +7.59%  '
$ws.Range("D45").Value = "'4.75"
$ws.Range("E45").Value = '  +14.40%  '
$ws.Range("D46").Value = "'0.0921"
$ws.Range("E46").Value = '  +6.74%  '
$ws.Range("D47").Value = "'259.19"
$ws.Range("E47").Value = '  +16.54%  '
$ws.Range("D48").Value = "'10.18"
$ws.Range("E48").Value = '  +0.69%  '
$ws.Range("D49").Value = "'0.0228"
$ws.Range("E49").Value = '  +5.70%  '
$ws.Range("D50").Value = '1.885.02'
$ws.Range("E50").Value = '  -2.68%  '
$ws.Range("D51").Value = "'17.47"
$ws.Range("E51").Value = '  +6.35%  '
